$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.486.70"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "2.985.25"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'381.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "

$ws.Range("D6").Value = "'103.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.67%  "

$ws.Range("E7").Value = "  +1.00%  "

$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("D10").Value = "'36.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("E11").Value = "  -0.90%  "

$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "3.456.35"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("D14").Value = "'18.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.78%  "

$ws.Range("E15").Value = "  +2.43%  "

$ws.Range("D16").Value = "2.978.27"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("D17").Value = "'11.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.31%  "

$ws.Range("D18").Value = "'0.996"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").Value = "51.498.27"
$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("E21").Value = "  +0.96%  "

$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").Value = "'70.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.15%  "

$ws.Range("D24").Value = "'267.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("E25").Value = "  +3.06%  "

$ws.Range("D26").Value = "'7.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.83%  "

$ws.Range("D27").Value = "'7.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.98%  "

$ws.Range("E28").Value = "  +2.55%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  +1.55%  "

$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").Value = "'10.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.41%  "

$ws.Range("D33").Value = "'34.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.86%  "

$ws.Range("D34").Value = "'51.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.10%  "

$ws.Range("D35").Value = "'2.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.25%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  +3.71%  "

$ws.Range("D39").Value = "'16.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.40%  "

$ws.Range("D40").Value = "'2.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.81%  "

$ws.Range("E41").Value = "  +1.00%  "

$ws.Range("D42").Value = "'1.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.51%  "

$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'3.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.86%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'123.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").Value = "'21.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "

$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("D47").Value = "'0.272"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").Value = "'2.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "

$ws.Range("D49").Value = "2.028.38"
$ws.Range("E49").Value = "  +1.88%  "

$ws.Range("D50").Value = "3.282.19"
$ws.Range("E50").Value = "  +1.42%  "

$ws.Range("D51").Value = "'0.0334"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.37%  "
